$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1508.125
$ws.Range("I2").Value = 1438.6923
$ws.Range("J2").Value = 1590.1818
$ws.Range("K2").Value = 1438.6923
$ws.Range("L2").Value = 1590.1818
$ws.Range("M2").Value = -1325.6923
$ws.Range("N2").Value = -1816.1818
$ws.Range("H133").Value = 90962.336
$ws.Range("J133").Value = 90962.336
$ws.Range("L133").Value = 90962.336
$ws.Range("N133").Value = -101082.336
$ws.Range("H134").Value = 134874.5
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140
$ws.Range("H136").Value = 131837.67
$ws.Range("J136").Value = 131837.67
$ws.Range("L136").Value = 131837.67
$ws.Range("N136").Value = -142037.67
$ws.Range("H137").Value = 2981.077
$ws.Range("I137").Value = 2965.5715
$ws.Range("K137").Value = 8896.7145
$ws.Range("M137").Value = -6346.7145
$ws.Range("H138").Value = 2103.0342
$ws.Range("I138").Value = 2011.9667
$ws.Range("J138").Value = 2150.138
$ws.Range("K138").Value = 6035.9001
$ws.Range("L138").Value = 6450.414
$ws.Range("M138").Value = -895.9000999999998
$ws.Range("N138").Value = -16730.414
$ws.Range("H140").Value = 210727
$ws.Range("J140").Value = 210727
$ws.Range("L140").Value = 210727
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1554.375
$ws.Range("I61").Value = 1471.7142
$ws.Range("J61").Value = 2133
$ws.Range("K61").Value = 1471.7142
$ws.Range("L61").Value = 2133
$ws.Range("M61").Value = -1259.7142
$ws.Range("N61").Value = -2557
$ws.Range("H74").Value = 1669.9445
$ws.Range("I74").Value = 1764.9333
$ws.Range("K74").Value = 1764.9333
$ws.Range("M74").Value = -890.9332999999999
$ws.Range("H77").Value = 1669.9445
$ws.Range("I77").Value = 1764.9333
$ws.Range("K77").Value = 8824.666499999999
$ws.Range("M77").Value = -4456.666499999999
$ws.Range("H102").Value = 1477.3158
$ws.Range("I102").Value = 1551.4615
$ws.Range("J102").Value = 1316.6666
$ws.Range("K102").Value = 1551.4615
$ws.Range("L102").Value = 1316.6666
$ws.Range("M102").Value = 70.53850000000011
$ws.Range("H132").Value = 2127.75
$ws.Range("I132").Value = 2185.375
$ws.Range("K132").Value = 6556.125
$ws.Range("M132").Value = -4026.125
$ws.Range("H136").Value = 1554.375
$ws.Range("I136").Value = 1471.7142
$ws.Range("J136").Value = 2133
$ws.Range("K136").Value = 4415.142599999999
$ws.Range("L136").Value = 6399
$ws.Range("M136").Value = -1865.142599999999
$ws.Range("N136").Value = -11499
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3708.25
$ws.Range("I86").Value = 3796.3333
$ws.Range("J86").Value = 3444
$ws.Range("K86").Value = 3796.3333
$ws.Range("L86").Value = 3444
$ws.Range("M86").Value = -2673.3333
$ws.Range("N86").Value = -5690
$ws.Range("H89").Value = 3708.25
$ws.Range("I89").Value = 3796.3333
$ws.Range("J89").Value = 3444
$ws.Range("K89").Value = 18981.6665
$ws.Range("L89").Value = 17220
$ws.Range("M89").Value = -13365.6665
$ws.Range("N89").Value = -28452
$ws.Range("H132").Value = 101488.664
$ws.Range("J132").Value = 101488.664
$ws.Range("L132").Value = 101488.664
$ws.Range("N132").Value = -111608.664
$ws.Range("H134").Value = 3188.3333
$ws.Range("I134").Value = 2849
$ws.Range("J134").Value = 3612.5
$ws.Range("K134").Value = 8547
$ws.Range("L134").Value = 10837.5
$ws.Range("M134").Value = -6012
$ws.Range("N134").Value = -15907.5
$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27223.291
$ws.Range("I31").Value = 29107.105
$ws.Range("J31").Value = 23012.412
$ws.Range("K31").Value = 29107.105
$ws.Range("L31").Value = 23012.412
$ws.Range("M31").Value = -28812.105
$ws.Range("N31").Value = -23602.412
$ws.Range("H34").Value = 27223.291
$ws.Range("I34").Value = 29107.105
$ws.Range("J34").Value = 23012.412
$ws.Range("K34").Value = 29107.105
$ws.Range("L34").Value = 23012.412
$ws.Range("M34").Value = -28905.105
$ws.Range("N34").Value = -23416.412
$ws.Range("H58").Value = 993.8214
$ws.Range("I58").Value = 1019.5
$ws.Range("K58").Value = 1019.5
$ws.Range("M58").Value = -816.5
$ws.Range("H132").Value = 2627.7917
$ws.Range("I132").Value = 2775.7727
$ws.Range("K132").Value = 8327.3181
$ws.Range("M132").Value = -5797.3181
$ws.Range("H134").Value = 2825.7925
$ws.Range("I134").Value = 2764.5
$ws.Range("J134").Value = 3014.3845
$ws.Range("K134").Value = 8293.5
$ws.Range("L134").Value = 9043.1535
$ws.Range("M134").Value = -5758.5
$ws.Range("N134").Value = -14113.1535
$ws.Range("H136").Value = 993.8214
$ws.Range("I136").Value = 1019.5
$ws.Range("K136").Value = 3058.5
$ws.Range("M136").Value = -508.5
$ws.Range("H140").Value = 99888
$ws.Range("J140").Value = 99888
$ws.Range("L140").Value = 99888
$ws.Range("N140").Value = -110248
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2674.0889
$ws.Range("I132").Value = 2054.457
$ws.Range("J132").Value = 4842.8
$ws.Range("K132").Value = 6163.370999999999
$ws.Range("L132").Value = 14528.4
$ws.Range("M132").Value = -3633.370999999999
$ws.Range("N132").Value = -19588.4
$ws.Range("H135").Value = 183333.33
$ws.Range("J135").Value = 183333.33
$ws.Range("L135").Value = 183333.33
$ws.Range("N135").Value = -193473.33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1558.8572
$ws.Range("I68").Value = 1558.8572
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1558.8572
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -809.8571999999999
$ws.Range("H71").Value = 1558.8572
$ws.Range("I71").Value = 1558.8572
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 7794.286
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4050.286
$ws.Range("H132").Value = 2061.61
$ws.Range("I132").Value = 1772.9814
$ws.Range("J132").Value = 5178.8
$ws.Range("K132").Value = 5318.9442
$ws.Range("L132").Value = 15536.4
$ws.Range("M132").Value = -2788.9442
$ws.Range("N132").Value = -20596.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2737.5962
$ws.Range("I132").Value = 2766.6458
$ws.Range("J132").Value = 2389
$ws.Range("K132").Value = 8299
$ws.Range("L132").Value = 7167
$ws.Range("M132").Value = -5769.937399999999
$ws.Range("N132").Value = -12227

# Additions (cells that did not exist before)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N140").Value = -221087
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N102").Value = -4560.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N137").Value = -110195

# Removals (cells cleared entirely)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
